$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: remove the two weekly data rows for
# 2024-03-16 (45361.99999999999 / 454) and 2024-03-23 (45368.99999999999 / 258).
# Deleting the sheet rows shifts everything below up by 2 (A1:B31 -> A1:B29).
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("6:7").Delete()

# "Monthly Trend" sheet: remove the monthly data row for
# 2024-04 (45382.99999999999 / 712).
# Deleting the sheet row shifts everything below up by 1 (A1:B12 -> A1:B11).
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows("4:4").Delete()
